$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 14
$ws.Range("H5").Value = 498.42856
$ws.Range("I5").Value = 531.6667
$ws.Range("J5").Value = 299
$ws.Range("K5").Value = 531.6667
$ws.Range("L5").Value = 299
$ws.Range("M5").Value = -416.6667
$ws.Range("N5").Value = -529
$ws.Range("H18").Value = 3199.3333
$ws.Range("I18").Value = 3199.3333
$ws.Range("K18").Value = 3199.3333
$ws.Range("M18").Value = -2915.3333
$ws.Range("H33").Value = 267.48148
$ws.Range("I33").Value = 275.35
$ws.Range("K33").Value = 275.35
$ws.Range("M33").Value = -46.35000000000002
$ws.Range("H80").Value = 419.42856
$ws.Range("J80").Value = 439.33334
$ws.Range("L80").Value = 1318.00002
$ws.Range("N80").Value = -3314.00002
$ws.Range("H83").Value = 419.42856
$ws.Range("J83").Value = 439.33334
$ws.Range("L83").Value = 3954.00006
$ws.Range("N83").Value = -13938.00006
$ws.Range("H103").Value = 450
$ws.Range("J103").Value = 325
$ws.Range("L103").Value = 975
$ws.Range("N103").Value = -2147
$ws.Range("H137").Value = 1503.3469
$ws.Range("I137").Value = 1377.8334
$ws.Range("K137").Value = 4133.5002
$ws.Range("M137").Value = -1583.5002
$ws.Range("H138").Value = 6669386.5
$ws.Range("I138").Value = 1441.3182
$ws.Range("J138").Value = 9437213
$ws.Range("K138").Value = 4323.9546
$ws.Range("L138").Value = 28311639
$ws.Range("M138").Value = 816.0454
$ws.Range("N138").Value = -28321919

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 643
$ws.Range("I4").Value = 180
$ws.Range("K4").Value = 180
$ws.Range("M4").Value = -64
$ws.Range("H5").Value = 51.565216
$ws.Range("I5").Value = 80.85714
$ws.Range("J5").Value = 38.75
$ws.Range("K5").Value = 80.85714
$ws.Range("L5").Value = 38.75
$ws.Range("M5").Value = 31.14286
$ws.Range("N5").Value = -262.75
$ws.Range("H32").Value = 6985.612
$ws.Range("I32").Value = 3392.8572
$ws.Range("J32").Value = 25276
$ws.Range("K32").Value = 3392.8572
$ws.Range("L32").Value = 25276
$ws.Range("M32").Value = -3105.8572
$ws.Range("N32").Value = -25850
$ws.Range("H74").Value = 6814.65
$ws.Range("I74").Value = 1905.1666
$ws.Range("J74").Value = 51000
$ws.Range("K74").Value = 1905.1666
$ws.Range("L74").Value = 51000
$ws.Range("M74").Value = -1031.1666
$ws.Range("N74").Value = -52748
$ws.Range("H77").Value = 6814.65
$ws.Range("I77").Value = 1905.1666
$ws.Range("J77").Value = 51000
$ws.Range("K77").Value = 9525.833000000001
$ws.Range("L77").Value = 255000
$ws.Range("M77").Value = -5157.833000000001
$ws.Range("N77").Value = -263736
$ws.Range("H132").Value = 3104.25
$ws.Range("I132").Value = 2796.6843
$ws.Range("K132").Value = 8390.052899999999
$ws.Range("M132").Value = -5860.052899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 51.565216
$ws.Range("I4").Value = 80.85714
$ws.Range("J4").Value = 38.75
$ws.Range("K4").Value = 80.85714
$ws.Range("L4").Value = 38.75
$ws.Range("M4").Value = 34.14286
$ws.Range("N4").Value = -268.75
$ws.Range("H134").Value = 2163.6072
$ws.Range("I134").Value = 2160.4075
$ws.Range("K134").Value = 6481.2225
$ws.Range("M134").Value = -3946.2225

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 262.4
$ws.Range("I7").Value = 247.57143
$ws.Range("J7").Value = 281.27274
$ws.Range("K7").Value = 247.57143
$ws.Range("L7").Value = 281.27274
$ws.Range("M7").Value = -134.57143
$ws.Range("N7").Value = -507.27274
$ws.Range("H10").Value = 804.3333
$ws.Range("I10").Value = 202.5
$ws.Range("K10").Value = 202.5
$ws.Range("M10").Value = -63.5
$ws.Range("H13").Value = 272.5
$ws.Range("I13").Value = 45
$ws.Range("K13").Value = 45
$ws.Range("M13").Value = 94
$ws.Range("H14").Value = 800
$ws.Range("I14").Value = 800
$ws.Range("K14").Value = 800
$ws.Range("M14").Value = -630
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H22").Value = 228.11111
$ws.Range("I22").Value = 228.11111
$ws.Range("K22").Value = 228.11111
$ws.Range("M22").Value = 121.88889
$ws.Range("H31").Value = 58387.39
$ws.Range("I31").Value = 73979.64
$ws.Range("J31").Value = 3814.5
$ws.Range("K31").Value = 73979.64
$ws.Range("L31").Value = 3814.5
$ws.Range("M31").Value = -73684.64
$ws.Range("N31").Value = -4404.5
$ws.Range("H34").Value = 58387.39
$ws.Range("I34").Value = 73979.64
$ws.Range("J34").Value = 3814.5
$ws.Range("K34").Value = 73979.64
$ws.Range("L34").Value = 3814.5
$ws.Range("M34").Value = -73777.64
$ws.Range("N34").Value = -4218.5
$ws.Range("H42").Value = 5899.5
$ws.Range("J42").Value = 5899.5
$ws.Range("L42").Value = 5899.5
$ws.Range("N42").Value = -7085.5
$ws.Range("H60").Value = 17822.111
$ws.Range("J60").Value = 17914.285
$ws.Range("L60").Value = 17914.285
$ws.Range("N60").Value = -18936.285
$ws.Range("H132").Value = 2964.3103
$ws.Range("I132").Value = 2900.4783
$ws.Range("J132").Value = 3209
$ws.Range("K132").Value = 8701.4349
$ws.Range("L132").Value = 9627
$ws.Range("M132").Value = -6171.4349
$ws.Range("N132").Value = -14687
$ws.Range("H134").Value = 18877.055
$ws.Range("I134").Value = 5220.375
$ws.Range("J134").Value = 106279.8
$ws.Range("K134").Value = 15661.125
$ws.Range("L134").Value = 318839.4
$ws.Range("M134").Value = -13126.125
$ws.Range("N134").Value = -323909.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1484.9
$ws.Range("I50").Value = 220
$ws.Range("J50").Value = 2027
$ws.Range("K50").Value = 660
$ws.Range("L50").Value = 6081
$ws.Range("M50").Value = -179
$ws.Range("N50").Value = -7043
$ws.Range("H53").Value = 1484.9
$ws.Range("I53").Value = 220
$ws.Range("J53").Value = 2027
$ws.Range("K53").Value = 660
$ws.Range("L53").Value = 6081
$ws.Range("M53").Value = -179
$ws.Range("N53").Value = -7043
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 11163.167
$ws.Range("I21").Value = 3993
$ws.Range("K21").Value = 3993
$ws.Range("M21").Value = -3820
$ws.Range("H30").Value = 11163.167
$ws.Range("I30").Value = 3993
$ws.Range("K30").Value = 3993
$ws.Range("M30").Value = -3888
$ws.Range("H70").Value = 3849.0667
$ws.Range("I70").Value = 3760.25
$ws.Range("K70").Value = 3760.25
$ws.Range("M70").Value = -3490.25
$ws.Range("H73").Value = 3849.0667
$ws.Range("I73").Value = 3760.25
$ws.Range("K73").Value = 3760.25
$ws.Range("M73").Value = -2824.25
$ws.Range("H122").Value = 2403.2703
$ws.Range("I122").Value = 2158.879
$ws.Range("J122").Value = 4419.5
$ws.Range("K122").Value = 6476.637
$ws.Range("L122").Value = 13258.5
$ws.Range("M122").Value = -4026.637
$ws.Range("N122").Value = -18158.5
$ws.Range("H132").Value = 3052.5652
$ws.Range("I132").Value = 2724.238
$ws.Range("K132").Value = 8172.714
$ws.Range("M132").Value = -5642.714
$ws.Range("H136").Value = 27713.562
$ws.Range("J136").Value = 27713.562
$ws.Range("L136").Value = 83140.686
$ws.Range("N136").Value = -88240.686

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3599.7
$ws.Range("I68").Value = 3599.7
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 3599.7
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2850.7
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 3599.7
$ws.Range("I71").Value = 3599.7
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 17998.5
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -14254.5
$ws.Range("N71").Value = 0
$ws.Range("H132").Value = 5292.5
$ws.Range("I132").Value = 4865.75
$ws.Range("J132").Value = 6999.5
$ws.Range("K132").Value = 14597.25
$ws.Range("L132").Value = 20998.5
$ws.Range("M132").Value = -12067.25
$ws.Range("N132").Value = -26058.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2355.9375
$ws.Range("I122").Value = 2044.3478
$ws.Range("K122").Value = 6133.0434
$ws.Range("M122").Value = -3683.0434
$ws.Range("H132").Value = 2018.2759
$ws.Range("I132").Value = 2021.875
$ws.Range("K132").Value = 6065.625
$ws.Range("M132").Value = -3535.625
